# Refresh the cryptocurrency price/volume snapshot (Price = column D,
# Volume(1h) = column E) with the values from the latest GitHub Actions run.
#
# Column D holds values such as "60.743.21" / "0.0000140" / "0.999" that
# look numeric but must stay literal text (leading/trailing zeros, the
# thousands-dot formatting, etc. all matter) -- so those cells are force-
# formatted as Text before the value is written, exactly like Excel marks a
# cell "number stored as text" when you need to keep the original notation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.743.21'
$ws.Range("E2").Value = '  -1.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.903.45'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '525.41'
$ws.Range("E5").Value = '  -2.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.42'
$ws.Range("E6").Value = '  -5.77%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.544'
$ws.Range("E8").Value = '  -4.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.910.55'
$ws.Range("E9").Value = '  -2.86%  '
$ws.Range("E10").Value = '  -4.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.06'
$ws.Range("E12").Value = '  -3.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.410.39'
$ws.Range("E13").Value = '  -2.92%  '
$ws.Range("E14").Value = '  +2.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.731.36'
$ws.Range("E15").Value = '  -1.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.50'
$ws.Range("E16").Value = '  -5.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.898.89'
$ws.Range("E17").Value = '  -3.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000140'
$ws.Range("E18").Value = '  -4.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.90'
$ws.Range("E19").Value = '  -4.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.57'
$ws.Range("E20").Value = '  -3.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '352.34'
$ws.Range("E21").Value = '  -7.42%  '
$ws.Range("E22").Value = '  -3.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("E24").Value = '  +1.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.98'
$ws.Range("E25").Value = '  -1.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.451'
$ws.Range("E26").Value = '  -4.24%  '
$ws.Range("E27").Value = '  -7.17%  '
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.83'
$ws.Range("E29").Value = '  -3.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0853'
$ws.Range("E30").Value = '  -9.40%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("E32").Value = '  -1.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.54'
$ws.Range("E33").Value = '  -4.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '153.09'
$ws.Range("E34").Value = '  -4.24%  '
$ws.Range("E35").Value = '  -4.53%  '
$ws.Range("E36").Value = '  -6.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.992'
$ws.Range("E37").Value = '  -7.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.19'
$ws.Range("E38").Value = '  -6.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.50'
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.46'
$ws.Range("E40").Value = '  -5.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.70'
$ws.Range("E41").Value = '  -5.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.283.76'
$ws.Range("E42").Value = '  -5.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.650'
$ws.Range("E43").Value = '  -3.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0580'
$ws.Range("E44").Value = '  -1.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.28'
$ws.Range("E45").Value = '  -8.01%  '
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.93'
$ws.Range("E47").Value = '  -4.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0237'
$ws.Range("E48").Value = '  -3.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.32'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0915'
$ws.Range("E50").Value = '  -4.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.38'
$ws.Range("E51").Value = '  -7.22%  '
